$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style/border/alignment) from the last existing data row (row 10)
# down into the two new rows, so the new rows end up with the same cell
# layout/styling as the rest of the table (year cell bold+bordered, data cells plain).
$ws.Range("A10:AQ10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:AQ10").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Row 11: 2021年 ----
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 841.3200000000001
$ws.Range("C11").Value = 302.55
$ws.Range("D11").Value = 35.96
$ws.Range("F11").Value = 412.31
$ws.Range("G11").Value = 2145.43
$ws.Range("H11").Value = 177.93
$ws.Range("I11").Value = 2212.79
$ws.Range("J11").Value = 90.84
$ws.Range("K11").Value = 22795.8
$ws.Range("L11").Value = 70.69
$ws.Range("M11").Value = 24.57
$ws.Range("N11").Value = 3.7
$ws.Range("O11").Value = 178.36
$ws.Range("P11").Value = 276.08
$ws.Range("Q11").Value = 52.04
$ws.Range("R11").Value = 15.06
$ws.Range("S11").Value = 509.53
$ws.Range("T11").Value = 134.64
$ws.Range("U11").Value = 3645.07
$ws.Range("W11").Value = 297.67
$ws.Range("X11").Value = 324.06
$ws.Range("Y11").Value = 145.33
$ws.Range("Z11").Value = 1166.32
$ws.Range("AA11").Value = 233.31
$ws.Range("AB11").Value = 339.95
$ws.Range("AC11").Value = 936.49
$ws.Range("AD11").Value = 235.09
$ws.Range("AE11").Value = 206.02
$ws.Range("AF11").Value = 3153.92
$ws.Range("AG11").Value = 1170.7
$ws.Range("AH11").Value = 416.1
$ws.Range("AI11").Value = 373.66
$ws.Range("AJ11").Value = 30.95
$ws.Range("AK11").Value = 572.89
$ws.Range("AL11").Value = 123.92
$ws.Range("AM11").Value = 836.88
$ws.Range("AN11").Value = 15.7
$ws.Range("AO11").Value = 676.29
$ws.Range("AP11").Value = 377.55
$ws.Range("AQ11").Value = 32.67

# ---- Row 12: 2022年 (only the total column K12 is populated) ----
$ws.Range("A12").Value = "2022年"
$ws.Range("K12").Value = 20040
